# Auto update Excel log
# Appends new sensor-event rows to the ALERTS, mmWave, Proximity, Camera and PIR
# sheets, mirroring a fresh batch of readings captured at 2026-01-30 ~15:21-15:23.

$wb = $excel.ActiveWorkbook

function Add-LogRows {
    param($SheetName, $StartRow, $Rows)

    $ws = $wb.Worksheets.Item($SheetName)
    $r = $StartRow
    foreach ($row in $Rows) {
        # Column A holds a literal date-like string (e.g. "2026-01-30"). Force
        # the cell to Text format first so Excel doesn't reinterpret it as a
        # real date serial number, keeping it consistent with the rest of the
        # log which stores everything as plain text.
        $ws.Cells.Item($r, 1).NumberFormat = "@"
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
        $ws.Cells.Item($r, 5).Value = $row[4]
        $ws.Cells.Item($r, 6).Value = $row[5]
        $r++
    }
}

# ALERTS: one new CRITICAL EMERGENCY / FALL_DETECTED row (row 12)
Add-LogRows "ALERTS" 12 @(
    ,@("2026-01-30", "15:21:56", "15:00", "Living Room", "CRITICAL EMERGENCY", "FALL_DETECTED")
)

# mmWave: seven new presence/fall rows (rows 80-86)
Add-LogRows "mmWave" 80 @(
    ,@("2026-01-30", "15:21:24", "15:00", "Living Room", "PRESENCE_DETECTED", "Active")
    ,@("2026-01-30", "15:21:34", "15:00", "Living Room", "PRESENCE_DETECTED", "Active")
    ,@("2026-01-30", "15:21:53", "15:00", "Living Room", "FALL_DETECTED", "EMERGENCY")
    ,@("2026-01-30", "15:22:49", "15:00", "Living Room", "PRESENCE_DETECTED", "Active")
    ,@("2026-01-30", "15:22:49", "15:00", "Living Room", "PRESENCE_DETECTED", "Active")
    ,@("2026-01-30", "15:22:49", "15:00", "Living Room", "PRESENCE_DETECTED", "Active")
    ,@("2026-01-30", "15:22:50", "15:00", "Living Room", "PRESENCE_DETECTED", "Active")
)

# Proximity: six new ENTER/EXIT rows (rows 32-37)
Add-LogRows "Proximity" 32 @(
    ,@("2026-01-30", "15:22:13", "15:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door")
    ,@("2026-01-30", "15:22:19", "15:00", "Living Room Main Door", "EXIT", "User EXITED Living Room Main Door")
    ,@("2026-01-30", "15:22:50", "15:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door")
    ,@("2026-01-30", "15:22:55", "15:00", "Living Room Main Door", "EXIT", "User EXITED Living Room Main Door")
    ,@("2026-01-30", "15:22:58", "15:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door")
    ,@("2026-01-30", "15:23:03", "15:00", "Living Room Main Door", "EXIT", "User EXITED Living Room Main Door")
)

# Camera: six new Image Captured rows (rows 31-36)
Add-LogRows "Camera" 31 @(
    ,@("2026-01-30", "15:22:13", "15:00", "Living Room Main Door", "Image Captured (ENTER)", "Active")
    ,@("2026-01-30", "15:22:19", "15:00", "Living Room Main Door", "Image Captured (EXIT)", "Active")
    ,@("2026-01-30", "15:22:50", "15:00", "Living Room Main Door", "Image Captured (ENTER)", "Active")
    ,@("2026-01-30", "15:22:55", "15:00", "Living Room Main Door", "Image Captured (EXIT)", "Active")
    ,@("2026-01-30", "15:22:57", "15:00", "Living Room Main Door", "Image Captured (ENTER)", "Active")
    ,@("2026-01-30", "15:23:03", "15:00", "Living Room Main Door", "Image Captured (EXIT)", "Active")
)

# PIR: three new RECOVERY_DETECTION rows (rows 23-25)
Add-LogRows "PIR" 23 @(
    ,@("2026-01-30", "15:21:23", "15:00", "Living Room", "RECOVERY_DETECTION", "Inactive")
    ,@("2026-01-30", "15:22:49", "15:00", "Living Room", "RECOVERY_DETECTION", "Inactive")
    ,@("2026-01-30", "15:22:49", "15:00", "Living Room", "RECOVERY_DETECTION", "Inactive")
)
